$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 425, shifting rows 425:503 down to 426:504.
$ws.Rows(425).Insert()

# Populate the newly inserted row 425 with the new daily price entry.
$ws.Cells.Item(425, 1).Value = 3
$ws.Cells.Item(425, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(425, 3).Value = "Coquimbo"
$ws.Cells.Item(425, 4).Value = 44782
$ws.Cells.Item(425, 5).Value = 5
$ws.Cells.Item(425, 6).Value = 100112003
$ws.Cells.Item(425, 7).Value = "Ajo"
$ws.Cells.Item(425, 8).Value = "Chino"
$ws.Cells.Item(425, 9).Value = "Primera"
$ws.Cells.Item(425, 10).Value = 78
$ws.Cells.Item(425, 11).Value = 24000
$ws.Cells.Item(425, 12).Value = 24500
$ws.Cells.Item(425, 13).Value = 24244
$ws.Cells.Item(425, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(425, 15).Value = "China"
$ws.Cells.Item(425, 16).Value = 2424
$ws.Cells.Item(425, 17).Value = 10
$ws.Cells.Item(425, 18).Value = "Hortaliza"
